$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cryptocurrency price (D) and volume change (E) columns.
# D-column values are stored as text in the workbook (e.g. "311.17" or
# European-grouped "45.103.13"). Plain numeric-looking strings get auto-
# converted to real numbers by Excel when assigned directly, so we prefix
# those with an apostrophe to force text entry, then clear the resulting
# cell format (quote-prefix) so the cell keeps the plain default style.

$c = $ws.Range("D2")
$c.Value = "'45.103.13"
$c.ClearFormats()
$ws.Range("E2").Value = "  +2.59%  "
$c = $ws.Range("D3")
$c.Value = "'2.370.69"
$c.ClearFormats()
$ws.Range("E3").Value = "  +1.18%  "
$ws.Range("E4").Value = "  -0.13%  "
$c = $ws.Range("D5")
$c.Value = "'311.17"
$c.ClearFormats()
$ws.Range("E5").Value = "  -0.72%  "
$c = $ws.Range("D6")
$c.Value = "'108.15"
$c.ClearFormats()
$ws.Range("E6").Value = "  -0.71%  "
$ws.Range("E7").Value = "  -0.16%  "
$c = $ws.Range("D9")
$c.Value = "'0.614"
$c.ClearFormats()
$ws.Range("E9").Value = "  -0.81%  "
$c = $ws.Range("D10")
$c.Value = "'40.87"
$c.ClearFormats()
$ws.Range("E10").Value = "  -0.73%  "
$c = $ws.Range("D11")
$c.Value = "'0.0917"
$c.ClearFormats()
$ws.Range("E11").Value = "  -0.35%  "
$c = $ws.Range("D12")
$c.Value = "'8.46"
$c.ClearFormats()
$ws.Range("E12").Value = "  -1.50%  "
$ws.Range("E13").Value = "  +1.42%  "
$c = $ws.Range("D14")
$c.Value = "'0.976"
$c.ClearFormats()
$ws.Range("E14").Value = "  -3.80%  "
$c = $ws.Range("D15")
$c.Value = "'2.734.33"
$c.ClearFormats()
$ws.Range("E15").Value = "  +1.38%  "
$c = $ws.Range("D16")
$c.Value = "'15.25"
$c.ClearFormats()
$ws.Range("E16").Value = "  -1.80%  "
$c = $ws.Range("D17")
$c.Value = "'2.387.12"
$c.ClearFormats()
$ws.Range("E17").Value = "  +2.44%  "
$c = $ws.Range("D18")
$c.Value = "'45.139.59"
$c.ClearFormats()
$ws.Range("E18").Value = "  +2.98%  "
$c = $ws.Range("D19")
$c.Value = "'14.60"
$c.ClearFormats()
$ws.Range("E19").Value = "  +11.96%  "
$c = $ws.Range("D20")
$c.Value = "'7.27"
$c.ClearFormats()
$ws.Range("E20").Value = "  -4.02%  "
$ws.Range("E21").Value = "  -0.71%  "
$c = $ws.Range("D22")
$c.Value = "'73.27"
$c.ClearFormats()
$ws.Range("E22").Value = "  -1.42%  "
$c = $ws.Range("D23")
$c.Value = "'3.49"
$c.ClearFormats()
$ws.Range("E23").Value = "  -0.28%  "
$c = $ws.Range("D24")
$c.Value = "'259.48"
$c.ClearFormats()
$ws.Range("E24").Value = "  -3.75%  "
$ws.Range("E25").Value = "  +0.98%  "
$ws.Range("E26").Value = "  +0.11%  "
$ws.Range("E27").Value = "  +0.02%  "
$c = $ws.Range("D28")
$c.Value = "'7.21"
$c.ClearFormats()
$ws.Range("E28").Value = "  -5.56%  "
$ws.Range("E29").Value = "  +1.78%  "
$c = $ws.Range("D30")
$c.Value = "'0.0964"
$c.ClearFormats()
$ws.Range("E30").Value = "  +8.53%  "
$c = $ws.Range("D31")
$c.Value = "'22.36"
$c.ClearFormats()
$ws.Range("E31").Value = "  -1.53%  "
$c = $ws.Range("D32")
$c.Value = "'37.37"
$c.ClearFormats()
$ws.Range("E32").Value = "  -3.85%  "
$c = $ws.Range("D33")
$c.Value = "'169.16"
$c.ClearFormats()
$ws.Range("E33").Value = "  +0.45%  "
$c = $ws.Range("D34")
$c.Value = "'2.95"
$c.ClearFormats()
$ws.Range("E34").Value = "  +6.29%  "
$ws.Range("E35").Value = "  -1.83%  "
$ws.Range("E36").Value = "  +2.57%  "
$c = $ws.Range("D37")
$c.Value = "'4.73"
$c.ClearFormats()
$ws.Range("E37").Value = "  -0.88%  "
$c = $ws.Range("D38")
$c.Value = "'3.92"
$c.ClearFormats()
$c = $ws.Range("D39")
$c.Value = "'2.92"
$c.ClearFormats()
$ws.Range("E39").Value = "  +0.93%  "
$c = $ws.Range("D40")
$c.Value = "'0.0353"
$c.ClearFormats()
$ws.Range("E40").Value = "  -3.67%  "
$ws.Range("E41").Value = "  +3.60%  "
$c = $ws.Range("D42")
$c.Value = "'99.70"
$c.ClearFormats()
$ws.Range("E42").Value = "  -4.46%  "
$c = $ws.Range("D43")
$c.Value = "'1.891.12"
$c.ClearFormats()
$ws.Range("E43").Value = "  +13.16%  "
$c = $ws.Range("D44")
$c.Value = "'69.90"
$c.ClearFormats()
$ws.Range("E44").Value = "  -2.63%  "
$ws.Range("E45").Value = "  -4.09%  "
$c = $ws.Range("D46")
$c.Value = "'12.81"
$c.ClearFormats()
$ws.Range("E46").Value = "  -4.01%  "
$ws.Range("E47").Value = "  -0.10%  "
$c = $ws.Range("D48")
$c.Value = "'81.28"
$c.ClearFormats()
$ws.Range("E48").Value = "  +5.57%  "
$c = $ws.Range("D49")
$c.Value = "'5.62"
$c.ClearFormats()
$ws.Range("E49").Value = "  +7.07%  "
$c = $ws.Range("D50")
$c.Value = "'112.34"
$c.ClearFormats()
$ws.Range("E50").Value = "  -1.61%  "
$c = $ws.Range("D51")
$c.Value = "'9.20"
$c.ClearFormats()
$ws.Range("E51").Value = "  +2.29%  "
